$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1991
$ws.Range("D2").Value = 7270
$ws.Range("E2").Value = 27.39
$ws.Range("F2").Value = 26.36
$ws.Range("G2").Value = 28.41
$ws.Range("C3").Value = 860
$ws.Range("D3").Value = 7270
$ws.Range("E3").Value = 11.83
$ws.Range("F3").Value = 11.09
$ws.Range("G3").Value = 12.57
$ws.Range("B4").Value = "Injuries & adverse effects"
$ws.Range("C4").Value = 769
$ws.Range("D4").Value = 7270
$ws.Range("F4").Value = 9.869999999999999
$ws.Range("G4").Value = 11.28
$ws.Range("B5").Value = "Symptom – Circulatory"
$ws.Range("C5").Value = 766
$ws.Range("D5").Value = 7270
$ws.Range("E5").Value = 10.54
$ws.Range("F5").Value = 9.83
$ws.Range("G5").Value = 11.24
$ws.Range("B6").Value = "Symptom – Nervous"
$ws.Range("C6").Value = 760
$ws.Range("D6").Value = 7270
$ws.Range("E6").Value = 10.45
$ws.Range("F6").Value = 9.75
$ws.Range("G6").Value = 11.16
$ws.Range("C7").Value = 500
$ws.Range("D7").Value = 7270
$ws.Range("E7").Value = 6.88
$ws.Range("F7").Value = 6.3
$ws.Range("G7").Value = 7.46
$ws.Range("B8").Value = "Diseases (patient-stated)"
$ws.Range("C8").Value = 470
$ws.Range("D8").Value = 7270
$ws.Range("E8").Value = 6.46
$ws.Range("F8").Value = 5.9
$ws.Range("G8").Value = 7.03
$ws.Range("B9").Value = "Symptom – Genitourinary"
$ws.Range("C9").Value = 424
$ws.Range("D9").Value = 7270
$ws.Range("E9").Value = 5.83
$ws.Range("F9").Value = 5.29
$ws.Range("G9").Value = 6.37
$ws.Range("C10").Value = 343
$ws.Range("D10").Value = 7270
$ws.Range("E10").Value = 4.72
$ws.Range("F10").Value = 4.23
$ws.Range("G10").Value = 5.21
$ws.Range("B11").Value = "Uncodable/Unknown"
$ws.Range("C11").Value = 194
$ws.Range("D11").Value = 7270
$ws.Range("E11").Value = 2.67
$ws.Range("F11").Value = 2.3
$ws.Range("G11").Value = 3.04
$ws.Range("B12").Value = "Symptom – Skin/Hair/Nails"
$ws.Range("C12").Value = 193
$ws.Range("D12").Value = 7270
$ws.Range("E12").Value = 2.65
$ws.Range("F12").Value = 2.29
$ws.Range("G12").Value = 3.02
$ws.Range("C17").Value = 107
$ws.Range("E17").Value = 5.4
$ws.Range("F17").Value = 4.4
$ws.Range("G17").Value = 6.39
$ws.Range("B23").Value = "Uncodable/Unknown"
$ws.Range("C23").Value = 21
$ws.Range("E23").Value = 1.06
$ws.Range("F23").Value = 0.61
$ws.Range("G23").Value = 1.51
$ws.Range("C24").Value = 462
$ws.Range("D24").Value = 1350
$ws.Range("E24").Value = 34.22
$ws.Range("F24").Value = 31.69
$ws.Range("G24").Value = 36.75
$ws.Range("B25").Value = "Symptom – Digestive"
$ws.Range("C25").Value = 158
$ws.Range("D25").Value = 1350
$ws.Range("E25").Value = 11.7
$ws.Range("F25").Value = 9.99
$ws.Range("G25").Value = 13.42
$ws.Range("B26").Value = "Symptom – Nervous"
$ws.Range("C26").Value = 152
$ws.Range("D26").Value = 1350
$ws.Range("E26").Value = 11.26
$ws.Range("F26").Value = 9.57
$ws.Range("G26").Value = 12.95
$ws.Range("C27").Value = 116
$ws.Range("D27").Value = 1350
$ws.Range("E27").Value = 8.59
$ws.Range("F27").Value = 7.1
$ws.Range("G27").Value = 10.09
$ws.Range("B28").Value = "Other"
$ws.Range("C28").Value = 101
$ws.Range("D28").Value = 1350
$ws.Range("E28").Value = 7.48
$ws.Range("F28").Value = 6.08
$ws.Range("G28").Value = 8.880000000000001
$ws.Range("B29").Value = "Injuries & adverse effects"
$ws.Range("C29").Value = 90
$ws.Range("D29").Value = 1350
$ws.Range("E29").Value = 6.67
$ws.Range("F29").Value = 5.34
$ws.Range("G29").Value = 8
$ws.Range("B30").Value = "Symptom – General"
$ws.Range("C30").Value = 74
$ws.Range("D30").Value = 1350
$ws.Range("E30").Value = 5.48
$ws.Range("F30").Value = 4.27
$ws.Range("G30").Value = 6.7
$ws.Range("B31").Value = "Symptom – Genitourinary"
$ws.Range("C31").Value = 68
$ws.Range("D31").Value = 1350
$ws.Range("E31").Value = 5.04
$ws.Range("F31").Value = 3.87
$ws.Range("G31").Value = 6.2
$ws.Range("B32").Value = "Diseases (patient-stated)"
$ws.Range("C32").Value = 66
$ws.Range("D32").Value = 1350
$ws.Range("E32").Value = 4.89
$ws.Range("F32").Value = 3.74
$ws.Range("G32").Value = 6.04
$ws.Range("C33").Value = 36
$ws.Range("D33").Value = 1350
$ws.Range("E33").Value = 2.67
$ws.Range("F33").Value = 1.81
$ws.Range("G33").Value = 3.53
$ws.Range("B34").Value = "Uncodable/Unknown"
$ws.Range("C34").Value = 27
$ws.Range("D34").Value = 1350
$ws.Range("E34").Value = 2
$ws.Range("F34").Value = 1.25
$ws.Range("G34").Value = 2.75
$ws.Range("C35").Value = 3386
$ws.Range("D35").Value = 11309
$ws.Range("E35").Value = 29.94
$ws.Range("F35").Value = 29.1
$ws.Range("G35").Value = 30.78
$ws.Range("B36").Value = "Symptom – Digestive"
$ws.Range("C36").Value = 1338
$ws.Range("D36").Value = 11309
$ws.Range("E36").Value = 11.83
$ws.Range("F36").Value = 11.24
$ws.Range("G36").Value = 12.43
$ws.Range("B37").Value = "Symptom – Nervous"
$ws.Range("C37").Value = 1297
$ws.Range("D37").Value = 11309
$ws.Range("E37").Value = 11.47
$ws.Range("F37").Value = 10.88
$ws.Range("G37").Value = 12.06
$ws.Range("C38").Value = 1141
$ws.Range("D38").Value = 11309
$ws.Range("E38").Value = 10.09
$ws.Range("F38").Value = 9.529999999999999
$ws.Range("G38").Value = 10.64
$ws.Range("C39").Value = 1017
$ws.Range("D39").Value = 11309
$ws.Range("E39").Value = 8.99
$ws.Range("F39").Value = 8.470000000000001
$ws.Range("G39").Value = 9.52
$ws.Range("C40").Value = 726
$ws.Range("D40").Value = 11309
$ws.Range("E40").Value = 6.42
$ws.Range("F40").Value = 5.97
$ws.Range("G40").Value = 6.87
$ws.Range("C41").Value = 686
$ws.Range("D41").Value = 11309
$ws.Range("E41").Value = 6.07
$ws.Range("F41").Value = 5.63
$ws.Range("G41").Value = 6.51
$ws.Range("B42").Value = "Symptom – Genitourinary"
$ws.Range("C42").Value = 645
$ws.Range("D42").Value = 11309
$ws.Range("E42").Value = 5.7
$ws.Range("F42").Value = 5.28
$ws.Range("G42").Value = 6.13
$ws.Range("B43").Value = "Symptom – General"
$ws.Range("C43").Value = 567
$ws.Range("D43").Value = 11309
$ws.Range("E43").Value = 5.01
$ws.Range("F43").Value = 4.61
$ws.Range("G43").Value = 5.42
$ws.Range("C44").Value = 272
$ws.Range("D44").Value = 11309
$ws.Range("E44").Value = 2.41
$ws.Range("F44").Value = 2.12
$ws.Range("G44").Value = 2.69
$ws.Range("B45").Value = "Uncodable/Unknown"
$ws.Range("C45").Value = 234
$ws.Range("D45").Value = 11309
$ws.Range("E45").Value = 2.07
$ws.Range("F45").Value = 1.81
$ws.Range("G45").Value = 2.33
$ws.Range("C46").Value = 2284
$ws.Range("D46").Value = 6244
$ws.Range("E46").Value = 36.58
$ws.Range("F46").Value = 35.38
$ws.Range("G46").Value = 37.77
$ws.Range("B47").Value = "Symptom – Nervous"
$ws.Range("C47").Value = 792
$ws.Range("D47").Value = 6244
$ws.Range("E47").Value = 12.68
$ws.Range("F47").Value = 11.86
$ws.Range("G47").Value = 13.51
$ws.Range("B48").Value = "Symptom – Digestive"
$ws.Range("C48").Value = 683
$ws.Range("D48").Value = 6244
$ws.Range("E48").Value = 10.94
$ws.Range("F48").Value = 10.16
$ws.Range("G48").Value = 11.71
$ws.Range("C49").Value = 538
$ws.Range("D49").Value = 6244
$ws.Range("E49").Value = 8.619999999999999
$ws.Range("F49").Value = 7.92
$ws.Range("G49").Value = 9.31
$ws.Range("B50").Value = "Injuries & adverse effects"
$ws.Range("C50").Value = 420
$ws.Range("D50").Value = 6244
$ws.Range("E50").Value = 6.73
$ws.Range("F50").Value = 6.11
$ws.Range("G50").Value = 7.35
$ws.Range("C51").Value = 361
$ws.Range("D51").Value = 6244
$ws.Range("E51").Value = 5.78
$ws.Range("F51").Value = 5.2
$ws.Range("G51").Value = 6.36
$ws.Range("B52").Value = "Symptom – Genitourinary"
$ws.Range("C52").Value = 331
$ws.Range("D52").Value = 6244
$ws.Range("E52").Value = 5.3
$ws.Range("F52").Value = 4.75
$ws.Range("G52").Value = 5.86
$ws.Range("B53").Value = "Symptom – General"
$ws.Range("C53").Value = 322
$ws.Range("D53").Value = 6244
$ws.Range("E53").Value = 5.16
$ws.Range("F53").Value = 4.61
$ws.Range("G53").Value = 5.71
$ws.Range("B54").Value = "Diseases (patient-stated)"
$ws.Range("C54").Value = 304
$ws.Range("D54").Value = 6244
$ws.Range("E54").Value = 4.87
$ws.Range("F54").Value = 4.33
$ws.Range("G54").Value = 5.4
$ws.Range("B55").Value = "Symptom – Skin/Hair/Nails"
$ws.Range("C55").Value = 126
$ws.Range("D55").Value = 6244
$ws.Range("E55").Value = 2.02
$ws.Range("F55").Value = 1.67
$ws.Range("G55").Value = 2.37
$ws.Range("B56").Value = "Uncodable/Unknown"
$ws.Range("C56").Value = 83
$ws.Range("D56").Value = 6244
$ws.Range("E56").Value = 1.33
$ws.Range("F56").Value = 1.05
$ws.Range("G56").Value = 1.61
